$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting pattern of row 94 into the two new rows (95 and 96).
# Only columns A:Y are copied because columns Z:AB in rows 95/96 already have the
# correct (unchanged) formatting in the source workbook.
$ws.Range("A94:Y94").Copy($ws.Range("A95"))
$ws.Range("A94:Y94").Copy($ws.Range("A96"))

# Fill in the two new validation rules (TPOD2190 and TPOD2200).
$ws.Range("B95").Value = "TPOD2190"
$ws.Range("C95").Value = "In het manifest-OW mag het objecttype Geometrie niet voorkomen."
$ws.Range("C96").Value = "In het manifest-OW mag een bestandsnaam niet eindigen op '.gml'"
$ws.Range("B96").Value = "TPOD2200"

# Remove the now-superfluous blank separator row, shifting the remaining rows up.
$ws.Rows(115).Delete()

# Put the selection on the first newly added row, as in the saved workbook.
$ws.Range("A95").Select()
